$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "Save" column (H) mirroring the formatting of the existing
# header column G (bold, bordered, centered header style).
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# Add the data value for the new column in row 2 (unstyled, like the
# other numeric data cells).
$ws.Range("H2").Value = 0
